$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "30.596.10"
$ws.Range('D3').Value = "1.879.02"
$ws.Range('E3').Value = "  -1.00%  "
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = "  -0.05%  "
$ws.Range('D5').Value = "'236.01"
$ws.Range('E5').Value = "  -4.01%  "
$ws.Range('D6').Value = "'1.001"
$ws.Range('E6').Value = "  -0.04%  "
$ws.Range('D7').Value = "'0.4872"
$ws.Range('E7').Value = "  -2.13%  "
$ws.Range('E8').Value = "  -2.33%  "
$ws.Range('D9').Value = "'0.06663"
$ws.Range('E9').Value = "  -2.41%  "
$ws.Range('D10').Value = "1.877.90"
$ws.Range('E10').Value = "  -1.00%  "
$ws.Range('D11').Value = "'16.62"
$ws.Range('E11').Value = "  -3.40%  "
$ws.Range('D12').Value = "'0.07234"
$ws.Range('E12').Value = "  -1.23%  "
$ws.Range('D13').Value = "'88.86"
$ws.Range('E13').Value = "  -2.85%  "
$ws.Range('D14').Value = "'4.998"
$ws.Range('E14').Value = "  -1.69%  "
$ws.Range('D15').Value = "'0.6520"
$ws.Range('E15').Value = "  -4.07%  "
$ws.Range('D16').Value = "30.537.38"
$ws.Range('E16').Value = "  -0.98%  "
$ws.Range('D17').Value = "'0.000007861"
$ws.Range('E17').Value = "  -2.30%  "
$ws.Range('E18').Value = "  -0.06%  "
$ws.Range('D19').Value = "'12.95"
$ws.Range('E19').Value = "  -3.31%  "
$ws.Range('D20').Value = "2.121.56"
$ws.Range('E20').Value = "  -1.38%  "
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = "  -0.06%  "
$ws.Range('D22').Value = "'4.707"
$ws.Range('D23').Value = "'192.21"
$ws.Range('E23').Value = "  +7.28%  "
$ws.Range('D24').Value = "'6.102"
$ws.Range('E24').Value = "  -0.06%  "
$ws.Range('D25').Value = "'9.278"
$ws.Range('E25').Value = "  -0.80%  "
$ws.Range('D26').Value = "'157.37"
$ws.Range('E26').Value = "  +1.26%  "
$ws.Range('D27').Value = "'18.37"
$ws.Range('E27').Value = "  -1.61%  "
$ws.Range('D28').Value = "'1.824"
$ws.Range('E28').Value = "  -6.24%  "
$ws.Range('D29').Value = "'1.404"
$ws.Range('E29').Value = "  +1.18%  "
$ws.Range('D30').Value = "'4.239"
$ws.Range('E30').Value = "  -2.82%  "
$ws.Range('D31').Value = "'0.08999"
$ws.Range('E31').Value = "  +0.45%  "
$ws.Range('D32').Value = "'3.920"
$ws.Range('E32').Value = "  -3.26%  "
$ws.Range('D33').Value = "'0.05125"
$ws.Range('E33').Value = "  -2.98%  "
$ws.Range('D34').Value = "'0.7215"
$ws.Range('E34').Value = "  -4.00%  "
$ws.Range('D35').Value = "'1.078"
$ws.Range('E35').Value = "  -5.84%  "
$ws.Range('D36').Value = "'2.693"
$ws.Range('E36').Value = "  -0.22%  "
$ws.Range('E37').Value = "  -4.84%  "
$ws.Range('E38').Value = "  -2.01%  "
$ws.Range('D39').Value = "'0.9165"
$ws.Range('E39').Value = "  -2.13%  "
$ws.Range('D40').Value = "'2.048"
$ws.Range('E40').Value = "  -6.41%  "
$ws.Range('D41').Value = "'0.4378"
$ws.Range('E41').Value = "  -0.15%  "
$ws.Range('D42').Value = "'104.61"
$ws.Range('E42').Value = "  -1.19%  "
$ws.Range('D43').Value = "'0.9960"
$ws.Range('E43').Value = "  -0.53%  "
$ws.Range('D44').Value = "'5.706"
$ws.Range('E44').Value = "  -2.79%  "
$ws.Range('D45').Value = "'0.1330"
$ws.Range('E45').Value = "  -3.69%  "
$ws.Range('D46').Value = "'7.345"
$ws.Range('E46').Value = "  -5.21%  "
$ws.Range('B47').Value = "Cronos"
$ws.Range('C47').Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('D47').Value = "'0.05822"
$ws.Range('E47').Value = "  -0.33%  "
$ws.Range('B48').Value = "Decentraland"
$ws.Range('C48').Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range('D48').Value = "'0.4005"
$ws.Range('E48').Value = "  +2.23%  "
$ws.Range('D49').Value = "'8.660"
$ws.Range('E49').Value = "  +0.63%  "
$ws.Range('D50').Value = "'1.403"
$ws.Range('E50').Value = "  +0.78%  "
$ws.Range('D51').Value = "'33.05"
$ws.Range('E51').Value = "  -1.72%  "
